# Revert "adding term 2.0.0"
# - restore Metadata sheet values (Version, Date, Contact)
# - restore the "Include from FSIII" concept value back to "A"
# - remove the extra "Include from FSIII 2" sheet that the original commit added

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: revert the three changed values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- Include from FSIII sheet: revert the concept value back to "A" ---
$include1 = $wb.Worksheets.Item("Include from FSIII")
$include1.Range("C2").Value = "A"

# --- Remove the sheet added by the reverted commit ---
$excel.DisplayAlerts = $false
$include2 = $wb.Worksheets.Item("Include from FSIII 2")
[void]$include2.Delete()

# Deleting a sheet can shift the active tab; restore the original
# active sheet ("Metadata", the first tab) to match the pre-edit selection.
$meta.Activate()
